# COTX added to ESAD
# Slide 5 of the FIDO wallet whitepaper: the "Decrypt and validate /
# authorization" textbox becomes "Decrypt and validate / user authorization",
# and the "Encrypt authorization" textbox becomes "Encrypt user authorization"
# (widened/shifted left to fit the new wording).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- "Decrypt and validate" / "authorization" textbox (TextBox 34, id 35) ---
# This shape holds a single paragraph: a run "Decrypt and validate", a
# soft line break, then a run "authorization". Replace just the second
# run's text (characters 22..34, 1-based) so the line break survives.
$shpDecrypt = $s.Shapes.Item(43)
$trDecrypt = $shpDecrypt.TextFrame.TextRange
$runDecrypt = $trDecrypt.Characters(22, 13)
$runDecrypt.Text = "user authorization"
# The shape auto-fits its height to the text (spAutoFit); restore the
# original box height that the extra wording shrank.
$shpDecrypt.Height = 39.805

# --- "Encrypt authorization" textbox (TextBox 4, id 5) ---
$shpEncrypt = $s.Shapes.Item(59)
$shpEncrypt.TextFrame.TextRange.Text = "Encrypt user authorization"

# Restore the autofit height PowerPoint shrank after the text change, then
# reposition/resize the box (values are in points; COM stores EMU = pt*12700).
$shpEncrypt.Height = 23.64874
$shpEncrypt.Left = 75.352
$shpEncrypt.Width = 188.3163
